$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data range as Text so numeric-looking strings (e.g. "1.00") are not
# auto-converted to numbers by Excel, while keeping the original (unstyled) appearance.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.841.94"
$ws.Range("E2").Value = "  -5.59%  "
$ws.Range("D3").Value = "2.542.34"
$ws.Range("E3").Value = "  -5.44%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "298.85"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("D6").Value = "91.86"
$ws.Range("E6").Value = "  -6.76%  "
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  -3.86%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  -5.83%  "
$ws.Range("D10").Value = "35.97"
$ws.Range("E10").Value = "  -6.65%  "
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  -5.18%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.116"
$ws.Range("E12").Value = "  +7.33%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "7.68"
$ws.Range("E13").Value = "  -5.15%  "
$ws.Range("D14").Value = "2.929.95"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").Value = "2.583.94"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "0.876"
$ws.Range("E16").Value = "  -6.14%  "
$ws.Range("E17").Value = "  -6.61%  "
$ws.Range("D18").Value = "42.829.60"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("D19").Value = "0.0₃0979"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  -5.41%  "
$ws.Range("E22").Value = "  -4.13%  "
$ws.Range("D23").Value = "254.99"
$ws.Range("E23").Value = "  -9.65%  "
$ws.Range("D24").Value = "2.92"
$ws.Range("E24").Value = "  -4.51%  "
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("D26").Value = "28.91"
$ws.Range("E26").Value = "  -6.75%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").Value = "  -4.04%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.11"
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "36.52"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("D31").Value = "6.04"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "152.51"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").Value = "3.37"
$ws.Range("E34").Value = "  -11.32%  "
$ws.Range("D35").Value = "2.13"
$ws.Range("E35").Value = "  -9.23%  "
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("E37").Value = "  -6.12%  "
$ws.Range("D38").Value = "17.61"
$ws.Range("E38").Value = "  +7.57%  "
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "23.14"
$ws.Range("E40").Value = "  -10.37%  "
$ws.Range("E41").Value = "  -6.25%  "
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").Value = "2.089.32"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "1.94"
$ws.Range("E45").Value = "  +20.97%  "
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "9.11"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "84.25"
$ws.Range("E48").Value = "  -9.74%  "
$ws.Range("D49").Value = "2.785.45"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").Value = "104.72"
$ws.Range("E50").Value = "  -6.46%  "
$ws.Range("E51").Value = "  -3.03%  "

# Restore the default (unstyled) cell style now that values are safely stored as text,
# so the saved XML has no stray s="n" attributes on these cells, matching the original.
$dataRange.Style = "Normal"

Write-Host "Applied all changes"